$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-09-03 20:41:19"
$ws.Columns.Item(5).ColumnWidth = 17.2159881591797
$ws.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet -----------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-09-03 20:41:15"
$ws.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- de-de sheet -----------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-09-03 20:41:19"
$ws.Columns.Item(3).ColumnWidth = 17.2159881591797
